# The sheet's demo table (header row A1:C1 styled bold/centered/bordered,
# plus the data rows A2:C8) is being removed entirely, leaving an empty
# sheet (used range collapses back to A1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Clear()
